$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record needs to be inserted above the current row 418.
# This pushes every existing row from 418..521 down by one (to 419..522),
# growing the used range from A1:R521 to A1:R522, exactly like the rest of
# the table already does week over week.
$ws.Rows.Item(418).Insert()

# Populate the freshly inserted row with the new record. Every field matches
# the record that used to sit in row 418 except for the reporting date (D)
# and the volume (J), which carry the new week's figures.
$ws.Range("A418").Value = 10
$ws.Range("B418").Value = "Vega Modelo de Temuco"
$ws.Range("C418").Value = "La Araucanía"
$ws.Range("D418").Value = 45135
$ws.Range("E418").Value = 9
$ws.Range("F418").Value = 100112017
$ws.Range("G418").Value = "Apio"
$ws.Range("H418").Value = "Americana (o)"
$ws.Range("I418").Value = "Primera"
$ws.Range("J418").Value = 85
$ws.Range("K418").Value = 8000
$ws.Range("L418").Value = 8000
$ws.Range("M418").Value = 8000
$ws.Range("N418").Value = "$/docena de matas"
$ws.Range("O418").Value = "Provincia del Elquí"
$ws.Range("P418").Value = 1333
$ws.Range("Q418").Value = 6
$ws.Range("R418").Value = "Hortaliza"
